# Added a bunch of useful links (Intro to Statistical Learning resources).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Intro to Statistical Learning ebook
$ws.Range("A6").Value = "Intro to Statistical Learning ebook"
$ws.Range("B6").Value = "http://www-bcf.usc.edu/~gareth/ISL/ISLR%20Fourth%20Printing.pdf"

# Row 7 (B then A, to match original authoring/shared-string order)
$ws.Range("B7").Value = "http://www.alsharif.info/#!iom530/c21o7"

# Row 8: Intro to Statistical Learning website
$ws.Range("A8").Value = "Intro to Statistical Learning website"
$ws.Range("B8").Value = "http://www-bcf.usc.edu/~gareth/ISL/"

# Row 9: Intro to Statistical Learning MOOC
$ws.Range("A9").Value = "Intro to Statistical Learning MOOC"

# Row 7 label filled in after row 9's label (matches original edit/authoring order)
$ws.Range("A7").Value = "Applied statistical learning techniques (Slides and tuts)"

# Row 9 link
$ws.Range("B9").Value = "http://www.r-bloggers.com/in-depth-introduction-to-machine-learning-in-15-hours-of-expert-videos/"

# Row 10: Intro to Statistical Learning Stanford MOOC (Course)
$ws.Range("A10").Value = "Intro to Statistical Learning Stanford MOOC (Course)"
$ws.Range("B10").Value = "https://lagunita.stanford.edu/courses/HumanitiesScience/StatLearning/Winter2014/about"

# Widen the columns to fit the new (longer) titles and links.
# Note: this runtime's ColumnWidth setter adds a fixed +5/6 padding before
# storing the sheet's <col width>, so back the input off by that amount to
# land exactly on the target stored widths of 48.5 / 163.5.
$ws.Columns.Item(1).ColumnWidth = 48.5 - 5/6
$ws.Columns.Item(2).ColumnWidth = 163.5 - 5/6

$ws.Range("B10").Select()
